$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "250.12"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.08"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.448"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05666"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8152"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9229"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1442"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07538"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03106"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03088"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09372"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.756"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001591"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04768"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005790"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006367"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005047"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001033"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001502"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.700"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.199"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002999"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04028"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006779"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1068"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002714"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008040"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005802"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4999"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002103"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01011"
